$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r=3; $r -le 19; $r++) {
    if ($ws.Cells.Item($r, 2).Value2 -eq "Rad Paul") {
        $ws.Cells.Item($r, 2).Value2 = "Paul Rad"
    }
}

$ws.Cells.Item(20, 2).Value2 = "Bogdan Rulea"
$ws.Cells.Item(20, 6).Value2 = $true

$ws.Cells.Item(21, 2).Value2 = "Ionut Porumb"
$ws.Cells.Item(21, 6).Value2 = $true

$rng = $ws.Range("B3:R21")
$rng.Sort($ws.Range("B3:B21"))

for ($r=3; $r -le 21; $r++) {
    $line = "$r : "
    for ($c=2; $c -le 17; $c++) {
        $line += "[" + $ws.Cells.Item($r,$c).Value2 + "]"
    }
    Write-Host $line
}
